$wb = $excel.ActiveWorkbook
$target = $wb.Worksheets.Item("ARpUIiRC")
$s = $wb.Worksheets.Add($target)
$s.Name = "Blank2"
Write-Output ($wb.Worksheets | ForEach-Object { $_.Name })
Write-Output $s.Range("B1").Value
